$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (incl. the date style used in column A) from the row
# immediately above (row 59) down into the new row 60, then overwrite
# the values with the new data point.
$ws.Range("A59:F59").Copy()
$ws.Range("A60:F60").PasteSpecial(-4122)

$ws.Range("A60").Value = 45597
$ws.Range("B60").Value = -0.483
$ws.Range("C60").Value = 0.376
$ws.Range("D60").Value = -0.788
$ws.Range("E60").Value = 0.394
$ws.Range("F60").Value = 1.04
